$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "First Move Piece", "Capture Target" and "King Check" columns
# (D, E, F). This shifts the old "Rating" column (G) left into D, and the
# table/autofilter/used-range shrink from A1:G36 down to A1:D36 along with
# it. Deleting all three in one call keeps it to a single shift operation.
$ws.Range("D1:F1").EntireColumn.Delete()

# New selection left after the edit.
$ws.Range("H8").Select()
